$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Activate()

$ws.Range("B3").Value = 0.5
$ws.Range("B4").Value = 0.5
$ws.Range("B5").Value = 0.5

$ws.Range("C4").Select()
